$d = $word.ActiveDocument
$s = $d.Styles.Add("X10", 2)
$s.Priority = 1
$s.Font.Bold = $true
$s.Font.Shading.BackgroundPatternColor = 65535
Write-Output "done"
